# Scheduled runner: refresh computed profit/loss figures on each class sheet.
# Values mirror an upstream recompute (commit: "chore: update Sheets via scheduled runner").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2273893.2
$ws.Range("J17").Value = 2674965.8
$ws.Range("L17").Value = 8024897.399999999
$ws.Range("N17").Value = -8025233.399999999
$ws.Range("H43").Value = 9750
$ws.Range("I43").Value = 5000
$ws.Range("K43").Value = 5000
$ws.Range("M43").Value = -4931
$ws.Range("H86").Value = 46079.7
$ws.Range("J86").Value = 73916.164
$ws.Range("L86").Value = 73916.164
$ws.Range("N86").Value = -76162.164
$ws.Range("H88").Value = 4999.5
$ws.Range("I88").Value = 4999
$ws.Range("K88").Value = 4999
$ws.Range("M88").Value = -4593
$ws.Range("H89").Value = 46079.7
$ws.Range("J89").Value = 73916.164
$ws.Range("L89").Value = 369580.82
$ws.Range("N89").Value = -380812.82
$ws.Range("H91").Value = 4999.5
$ws.Range("I91").Value = 4999
$ws.Range("K91").Value = 4999
$ws.Range("M91").Value = -3595
$ws.Range("H106").Value = 5617114
$ws.Range("I106").Value = 7719973
$ws.Range("K106").Value = 7719973
$ws.Range("M106").Value = -7719342
$ws.Range("H107").Value = 6254.25
$ws.Range("I107").Value = 8328.625
$ws.Range("K107").Value = 8328.625
$ws.Range("M107").Value = -6408.625
$ws.Range("H125").Value = 4582.3335
$ws.Range("I125").Value = 5999.8887
$ws.Range("J125").Value = 3519.1667
$ws.Range("K125").Value = 53998.99830000001
$ws.Range("L125").Value = 31672.5003
$ws.Range("M125").Value = -51538.99830000001
$ws.Range("N125").Value = -36592.5003
$ws.Range("H131").Value = 6478
$ws.Range("I131").Value = 2956
$ws.Range("K131").Value = 8868
$ws.Range("M131").Value = -3828
$ws.Range("H138").Value = 378058.5
$ws.Range("I138").Value = 507225.16
$ws.Range("J138").Value = 4910.3335
$ws.Range("K138").Value = 1521675.48
$ws.Range("L138").Value = 14731.0005
$ws.Range("M138").Value = -1516535.48
$ws.Range("N138").Value = -25011.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3005.9312
$ws.Range("I32").Value = 3005.9312
$ws.Range("K32").Value = 3005.9312
$ws.Range("M32").Value = -2718.9312
$ws.Range("H61").Value = 5738.857
$ws.Range("I61").Value = 5770.8
$ws.Range("J61").Value = 5100
$ws.Range("K61").Value = 5770.8
$ws.Range("L61").Value = 5100
$ws.Range("M61").Value = -5558.8
$ws.Range("N61").Value = -5524
$ws.Range("H69").Value = 223555
$ws.Range("J69").Value = 223555
$ws.Range("L69").Value = 223555
$ws.Range("N69").Value = -225053
$ws.Range("H72").Value = 223555
$ws.Range("J72").Value = 223555
$ws.Range("L72").Value = 670665
$ws.Range("N72").Value = -678153
$ws.Range("H74").Value = 3283.4666
$ws.Range("I74").Value = 2020.9166
$ws.Range("J74").Value = 8333.666999999999
$ws.Range("K74").Value = 2020.9166
$ws.Range("L74").Value = 8333.666999999999
$ws.Range("M74").Value = -1146.9166
$ws.Range("N74").Value = -10081.667
$ws.Range("H77").Value = 3283.4666
$ws.Range("I77").Value = 2020.9166
$ws.Range("J77").Value = 8333.666999999999
$ws.Range("K77").Value = 10104.583
$ws.Range("L77").Value = 41668.335
$ws.Range("M77").Value = -5736.583000000001
$ws.Range("N77").Value = -50404.335
$ws.Range("H102").Value = 15788.667
$ws.Range("I102").Value = 25044.445
$ws.Range("J102").Value = 6532.8887
$ws.Range("K102").Value = 25044.445
$ws.Range("L102").Value = 6532.8887
$ws.Range("M102").Value = -23422.445
$ws.Range("N102").Value = -9776.8887
$ws.Range("H122").Value = 405018.6
$ws.Range("J122").Value = 1173541.8
$ws.Range("L122").Value = 3520625.4
$ws.Range("N122").Value = -3525525.4
$ws.Range("H136").Value = 5738.857
$ws.Range("I136").Value = 5770.8
$ws.Range("J136").Value = 5100
$ws.Range("K136").Value = 17312.4
$ws.Range("L136").Value = 15300
$ws.Range("M136").Value = -14762.4
$ws.Range("N136").Value = -20400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 12505917
$ws.Range("I7").Value = 6375.75
$ws.Range("J7").Value = 37505000
$ws.Range("K7").Value = 6375.75
$ws.Range("L7").Value = 37505000
$ws.Range("M7").Value = -6262.75
$ws.Range("N7").Value = -37505226
$ws.Range("H25").Value = 6136.6665
$ws.Range("I25").Value = 3608.6667
$ws.Range("K25").Value = 3608.6667
$ws.Range("M25").Value = -3373.6667
$ws.Range("H39").Value = 12559
$ws.Range("I39").Value = 12559
$ws.Range("K39").Value = 12559
$ws.Range("M39").Value = -12170
$ws.Range("H99").Value = 20729.191
$ws.Range("I99").Value = 34880.273
$ws.Range("J99").Value = 5163
$ws.Range("K99").Value = 34880.273
$ws.Range("L99").Value = 5163
$ws.Range("M99").Value = -33382.273
$ws.Range("N99").Value = -8159

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 304628.56
$ws.Range("J105").Value = 4466.6665
$ws.Range("L105").Value = 4466.6665
$ws.Range("N105").Value = -7960.6665
$ws.Range("H118").Value = 74999
$ws.Range("J118").Value = 74999
$ws.Range("L118").Value = 74999
$ws.Range("N118").Value = -78313

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1327.5454
$ws.Range("J34").Value = 1299.6666
$ws.Range("L34").Value = 3898.9998
$ws.Range("N34").Value = -4066.9998
$ws.Range("H39").Value = 840
$ws.Range("H55").Value = 5803.231
$ws.Range("J55").Value = 6914.8
$ws.Range("L55").Value = 20744.4
$ws.Range("N55").Value = -21098.4
$ws.Range("H94").Value = 2674.875
$ws.Range("J94").Value = 3325
$ws.Range("L94").Value = 9975
$ws.Range("N94").Value = -11327
$ws.Range("H137").Value = 4488.1113
$ws.Range("J137").Value = 10740.5
$ws.Range("L137").Value = 32221.5
$ws.Range("N137").Value = -42421.5
$ws.Range("H139").Value = 2354532.8
$ws.Range("I139").Value = 2858218.2
$ws.Range("K139").Value = 8574654.600000001
$ws.Range("M139").Value = -8569514.600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 18060
$ws.Range("I80").Value = 26400
$ws.Range("J80").Value = 12500
$ws.Range("K80").Value = 26400
$ws.Range("L80").Value = 12500
$ws.Range("M80").Value = -25402
$ws.Range("N80").Value = -14496
$ws.Range("H83").Value = 18060
$ws.Range("I83").Value = 26400
$ws.Range("J83").Value = 12500
$ws.Range("K83").Value = 132000
$ws.Range("L83").Value = 62500
$ws.Range("M83").Value = -127008
$ws.Range("N83").Value = -72484
$ws.Range("H126").Value = 9158.1
$ws.Range("I126").Value = 10622.2
$ws.Range("J126").Value = 7694
$ws.Range("K126").Value = 31866.6
$ws.Range("L126").Value = 23082
$ws.Range("M126").Value = -29396.6
$ws.Range("N126").Value = -28022
$ws.Range("H132").Value = 1371.0834
$ws.Range("I132").Value = 1359.3636
$ws.Range("K132").Value = 4078.0908
$ws.Range("M132").Value = -1548.0908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 8000
$ws.Range("I4").Value = 8000
$ws.Range("K4").Value = 8000
$ws.Range("M4").Value = -7887
$ws.Range("H25").Value = 18000
$ws.Range("I25").Value = 22000
$ws.Range("J25").Value = 10000
$ws.Range("K25").Value = 22000
$ws.Range("L25").Value = 10000
$ws.Range("M25").Value = -21770
$ws.Range("N25").Value = -10460
$ws.Range("H28").Value = 8000
$ws.Range("I28").Value = 8000
$ws.Range("K28").Value = 8000
$ws.Range("M28").Value = -7768
$ws.Range("H37").Value = 8000
$ws.Range("I37").Value = 8000
$ws.Range("K37").Value = 8000
$ws.Range("M37").Value = -7893
$ws.Range("H82").Value = 3759.5
$ws.Range("I82").Value = 4445
$ws.Range("J82").Value = 2799.8
$ws.Range("K82").Value = 4445
$ws.Range("L82").Value = 2799.8
$ws.Range("M82").Value = -4084
$ws.Range("N82").Value = -3521.8
$ws.Range("H85").Value = 3759.5
$ws.Range("I85").Value = 4445
$ws.Range("J85").Value = 2799.8
$ws.Range("K85").Value = 4445
$ws.Range("L85").Value = 2799.8
$ws.Range("M85").Value = -3197
$ws.Range("N85").Value = -5295.8
$ws.Range("H100").Value = 7597.8
$ws.Range("J100").Value = 13995
$ws.Range("L100").Value = 13995
$ws.Range("N100").Value = -15077
$ws.Range("H136").Value = 3374.652
$ws.Range("I136").Value = 2690.85
$ws.Range("K136").Value = 8072.549999999999
$ws.Range("M136").Value = -5522.549999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H126").Value = 60999.855
$ws.Range("I126").Value = 81400.8
$ws.Range("K126").Value = 244202.4
$ws.Range("M126").Value = -241732.4
$ws.Range("H139").Value = 59997
$ws.Range("J139").Value = 59997
$ws.Range("L139").Value = 59997
$ws.Range("N139").Value = -70277
